# Update the "Correspond Handoff Datetime" (D3) and "Correspond Handback DateTime" (G3)
# timestamps on the zh-cn and de-de report sheets to reflect a freshly regenerated
# handback report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D3").Value = "2016-01-13 11:39:26"
$wsZh.Range("G3").Value = "2016-01-13 11:40:45"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D3").Value = "2016-01-13 11:39:49"
$wsDe.Range("G3").Value = "2016-01-13 11:41:22"
